$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.341.07"
$ws.Range("E2").Value = "  +5.61%  "
$ws.Range("D3").Value = "1.921.60"
$ws.Range("E3").Value = "  +6.26%  "
$ws.Range("D4").Value = "'0.9978"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'254.70"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").Value = "'0.9981"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.5177"
$ws.Range("E7").Value = "  +4.35%  "
$ws.Range("D8").Value = "'46.28"
$ws.Range("E8").Value = "  +7.68%  "
$ws.Range("D9").Value = "'0.2982"
$ws.Range("E9").Value = "  +6.72%  "
$ws.Range("D10").Value = "'0.06828"
$ws.Range("E10").Value = "  +6.98%  "
$ws.Range("D11").Value = "1.917.45"
$ws.Range("E11").Value = "  +6.01%  "
$ws.Range("D12").Value = "'17.58"
$ws.Range("E12").Value = "  +4.79%  "
$ws.Range("D13").Value = "'0.07332"
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("D14").Value = "'0.6939"
$ws.Range("E14").Value = "  +7.41%  "
$ws.Range("D15").Value = "'88.38"
$ws.Range("E15").Value = "  +7.82%  "
$ws.Range("D16").Value = "'4.958"
$ws.Range("E16").Value = "  +5.57%  "
$ws.Range("D17").Value = "30.341.27"
$ws.Range("E17").Value = "  +5.65%  "
$ws.Range("D18").Value = "'0.000007969"
$ws.Range("E18").Value = "  +8.42%  "
$ws.Range("D19").Value = "'0.9961"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").Value = "'13.15"
$ws.Range("E20").Value = "  +7.29%  "
$ws.Range("D21").Value = "2.162.83"
$ws.Range("E21").Value = "  +6.02%  "
$ws.Range("D22").Value = "'0.9975"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "'4.889"
$ws.Range("E23").Value = "  +6.11%  "
$ws.Range("D24").Value = "'5.781"
$ws.Range("E24").Value = "  +8.86%  "
$ws.Range("D25").Value = "'9.251"
$ws.Range("E25").Value = "  +4.35%  "
$ws.Range("D26").Value = "'140.46"
$ws.Range("E26").Value = "  +26.11%  "
$ws.Range("D27").Value = "'147.46"
$ws.Range("E27").Value = "  +3.40%  "
$ws.Range("D28").Value = "'17.42"
$ws.Range("E28").Value = "  +8.85%  "
$ws.Range("D29").Value = "'2.028"
$ws.Range("E29").Value = "  +7.91%  "
$ws.Range("D30").Value = "'1.379"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("D31").Value = "'4.309"
$ws.Range("E31").Value = "  +3.04%  "
$ws.Range("D32").Value = "'0.08884"
$ws.Range("E32").Value = "  +6.30%  "
$ws.Range("D33").Value = "'4.050"
$ws.Range("E33").Value = "  +5.52%  "
$ws.Range("D34").Value = "'0.05153"
$ws.Range("E34").Value = "  +4.23%  "
$ws.Range("D35").Value = "'1.168"
$ws.Range("E35").Value = "  +6.95%  "
$ws.Range("D36").Value = "'0.7236"
$ws.Range("E36").Value = "  +7.80%  "
$ws.Range("E37").Value = "  +0.88%  "
$ws.Range("D38").Value = "'2.871"
$ws.Range("E38").Value = "  +8.83%  "
$ws.Range("E39").Value = "  +8.12%  "
$ws.Range("D40").Value = "'0.9742"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").Value = "'0.01707"
$ws.Range("E41").Value = "  +6.97%  "
$ws.Range("D42").Value = "'6.179"
$ws.Range("E42").Value = "  +4.37%  "
$ws.Range("D43").Value = "'0.4355"
$ws.Range("E43").Value = "  +5.85%  "
$ws.Range("D44").Value = "'106.36"
$ws.Range("E44").Value = "  +5.20%  "
$ws.Range("D45").Value = "'0.9990"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "'7.733"
$ws.Range("E46").Value = "  +6.90%  "
$ws.Range("D47").Value = "'0.1284"
$ws.Range("E47").Value = "  +4.67%  "
$ws.Range("D48").Value = "'0.05739"
$ws.Range("E48").Value = "  +4.46%  "
$ws.Range("D49").Value = "'8.587"
$ws.Range("E49").Value = "  +4.62%  "
$ws.Range("D50").Value = "'33.32"
$ws.Range("E50").Value = "  +6.34%  "
$ws.Range("D51").Value = "'0.3867"
$ws.Range("E51").Value = "  +7.26%  "
